$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Run ID value (B2). "VotingEnsemble" (B1) stays the same.
$ws.Range("B2").Value = "217f4a22-5141-4c51-9f1c-543416fef430_40"

# Update metric values (B3:B23)
$ws.Range("B3").Value = 0.82857000000000003
$ws.Range("B4").Value = 0.625
$ws.Range("B5").Value = 0.80815999999999999
$ws.Range("B6").Value = 0.625
$ws.Range("B7").Value = 0.70647000000000004
$ws.Range("B8").Value = 0.80506999999999995
$ws.Range("B9").Value = 0.80815999999999999
$ws.Range("B10").Value = 0.69167000000000001
$ws.Range("B11").Value = 0.66795000000000004
$ws.Range("B12").Value = 0.82857000000000003
$ws.Range("B13").Value = 0.79981999999999998
$ws.Range("B14").Value = 0.55188000000000004
$ws.Range("B15").Value = 0.39288000000000001
$ws.Range("B16").Value = 0.38333
$ws.Range("B17").Value = 0.67310000000000003
$ws.Range("B18").Value = 0.82857000000000003
$ws.Range("B19").Value = 0.80245
$ws.Range("B20").Value = 0.69167000000000001
$ws.Range("B21").Value = 0.82857000000000003
$ws.Range("B22").Value = 0.82857000000000003
$ws.Range("B23").Value = 0.86919000000000002

# Update the selected cell to match the target workbook's cursor position
$ws.Range("F7").Select()
